$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D (Price) updates.
# Values that still parse as a plain number (e.g. "0.998") must be forced to
# text with a leading apostrophe, exactly like a user typing '0.998 into a
# General-formatted cell, so the stored type stays text -- matching the
# original workbook's inlineStr ("41.528.25" style multi-dot numbers already
# stay textual on their own since they don't parse as numbers).
$dValues = @{
    2  = "41.528.55"
    3  = "2.485.84"
    4  = "'0.998"
    5  = "'313.36"
    6  = "'93.25"
    10 = "'32.60"
    13 = "2.867.59"
    15 = "'15.68"
    16 = "2.518.96"
    17 = "'0.752"
    18 = "41.552.70"
    19 = "'6.34"
    21 = "'70.83"
    22 = "'11.21"
    23 = "'236.20"
    27 = "'24.89"
    30 = "'36.38"
    31 = "'157.26"
    34 = "'18.11"
    35 = "'0.0755"
    38 = "'0.105"
    41 = "'4.11"
    43 = "'19.79"
    44 = "1.964.19"
    47 = "'8.87"
    48 = "2.726.77"
    49 = "'96.55"
    50 = "'67.53"
    51 = "'73.53"
}

# Column E (Volume(1h)) updates -- every data row (2..51) changes.
$eValues = @{
    2  = "  +0.37%  "
    3  = "  +0.77%  "
    4  = "  -0.22%  "
    5  = "  +0.48%  "
    6  = "  -1.22%  "
    7  = "  -0.90%  "
    8  = "  -0.19%  "
    9  = "  -0.98%  "
    10 = "  -3.34%  "
    11 = "  +0.48%  "
    12 = "  +1.94%  "
    13 = "  +0.74%  "
    14 = "  -2.19%  "
    15 = "  +7.08%  "
    16 = "  +0.73%  "
    17 = "  -4.60%  "
    18 = "  +0.58%  "
    19 = "  -0.01%  "
    20 = "  +1.17%  "
    21 = "  +4.27%  "
    22 = "  -2.59%  "
    23 = "  -0.19%  "
    24 = "  -2.87%  "
    25 = "  -0.17%  "
    26 = "  -1.40%  "
    27 = "  +1.80%  "
    28 = "  -0.03%  "
    29 = "  -0.56%  "
    30 = "  +0.48%  "
    31 = "  +2.57%  "
    32 = "  -2.57%  "
    33 = "  -1.27%  "
    34 = "  +5.62%  "
    35 = "  -0.10%  "
    36 = "  -5.36%  "
    37 = "  -2.24%  "
    38 = "  +2.29%  "
    39 = "  -3.12%  "
    40 = "  -0.22%  "
    41 = "  -3.77%  "
    42 = "  -0.24%  "
    43 = "  -7.07%  "
    44 = "  +0.02%  "
    45 = "  -0.22%  "
    46 = "  -3.54%  "
    47 = "  +2.13%  "
    48 = "  +0.72%  "
    49 = "  -0.82%  "
    50 = "  -3.39%  "
    51 = "  -3.51%  "
}

foreach ($row in $dValues.Keys) {
    $ws.Cells.Item($row, 4).Value = $dValues[$row]
}

foreach ($row in $eValues.Keys) {
    $ws.Cells.Item($row, 5).Value = $eValues[$row]
}
